# Refresh the "as_of_utc" timestamp stamped in column AA (rows 2-26) of the
# two stats sheets ("Главные" and "Линейные") from 2025-12-09 03:02:57 to
# 2025-12-09 07:02:50. The glossary sheet ("Глоссарий") has no AA column and
# is left untouched.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-12-09 03:02:57"
$newTimestamp = "2025-12-09 07:02:50"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $range = $ws.Range("AA2:AA26")
    $range.Value = $newTimestamp
}

Write-Output "Updated as_of_utc timestamps from $oldTimestamp to $newTimestamp on sheets: $sheetNames"
